$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number must be forced to Text format
# first (mirrors typing into a Text-formatted cell in real Excel), otherwise Excel
# auto-converts them to numeric values, same as it would for a live user editing session.

$ws.Range("D2").Value2 = "69.232.19"
$ws.Range("E2").Value2 = "  +2.02%  "
$ws.Range("D3").Value2 = "3.730.94"
$ws.Range("E3").Value2 = "  +0.42%  "
$ws.Range("E4").Value2 = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "612.32"
$ws.Range("E5").Value2 = "  +6.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "192.87"
$ws.Range("E6").Value2 = "  +10.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.640"
$ws.Range("E7").Value2 = "  +1.89%  "
$ws.Range("E8").Value2 = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.730"
$ws.Range("E9").Value2 = "  +3.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.163"
$ws.Range("E10").Value2 = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "60.26"
$ws.Range("E11").Value2 = "  +14.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.0000292"
$ws.Range("E12").Value2 = "  -1.64%  "
$ws.Range("E13").Value2 = "  +0.30%  "
$ws.Range("D14").Value2 = "4.323.39"
$ws.Range("E14").Value2 = "  -0.27%  "
$ws.Range("D15").Value2 = "3.727.53"
$ws.Range("E15").Value2 = "  -0.73%  "
$ws.Range("E16").Value2 = "  +2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "19.58"
$ws.Range("E17").Value2 = "  +1.08%  "
$ws.Range("E18").Value2 = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "12.99"
$ws.Range("E19").Value2 = "  +0.52%  "
$ws.Range("D20").Value2 = "69.069.91"
$ws.Range("E20").Value2 = "  +1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "413.64"
$ws.Range("E21").Value2 = "  +1.63%  "
$ws.Range("E22").Value2 = "  +2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "90.25"
$ws.Range("E23").Value2 = "  +2.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "3.10"
$ws.Range("E24").Value2 = "  +1.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "11.42"
$ws.Range("E25").Value2 = "  +7.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "12.99"
$ws.Range("E26").Value2 = "  +1.94%  "
$ws.Range("E27").Value2 = "  +1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "3.81"
$ws.Range("E28").Value2 = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "9.82"
$ws.Range("E29").Value2 = "  +3.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "33.06"
$ws.Range("E30").Value2 = "  +0.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "7.83"
$ws.Range("E31").Value2 = "  +2.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "12.83"
$ws.Range("E32").Value2 = "  +1.70%  "
$ws.Range("B33").Value2 = "Bittensor"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "650.81"
$ws.Range("E33").Value2 = "  +7.75%  "
$ws.Range("B34").Value2 = "Hedera"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.123"
$ws.Range("E34").Value2 = "  +5.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "46.00"
$ws.Range("E35").Value2 = "  +6.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "67.19"
$ws.Range("E36").Value2 = "  +3.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.420"
$ws.Range("E37").Value2 = "  +5.29%  "
$ws.Range("D38").Value2 = "0.0₃0837"
$ws.Range("E38").Value2 = "  -7.32%  "
$ws.Range("E39").Value2 = "  -0.07%  "
$ws.Range("E40").Value2 = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.142"
$ws.Range("E41").Value2 = "  +4.34%  "
$ws.Range("E42").Value2 = "  +1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.0451"
$ws.Range("E43").Value2 = "  +2.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "2.66"
$ws.Range("E44").Value2 = "  +3.47%  "
$ws.Range("D45").Value2 = "2.911.38"
$ws.Range("E45").Value2 = "  +5.92%  "
$ws.Range("E46").Value2 = "  +4.42%  "
$ws.Range("E47").Value2 = "  -0.33%  "
$ws.Range("E48").Value2 = "  +1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "3.10"
$ws.Range("E49").Value2 = "  -0.99%  "
$ws.Range("B50").Value2 = "dogwifhat"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "2.64"
$ws.Range("E50").Value2 = "  -12.03%  "
$ws.Range("B51").Value2 = "Monero"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "143.28"
$ws.Range("E51").Value2 = "  -0.05%  "
